# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume/listing updates described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.561.01'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.318.27'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.30'
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.66'
$ws.Range("E6").Value = '  -2.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.536'
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.342.08'
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.30'
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.737.20'
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("E15").Value = '  -4.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.436.04'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.343.82'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.20'
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("E21").Value = '  +2.31%  '
$ws.Range("E22").Value = '  -1.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.16'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.994'
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.34'
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.14'
$ws.Range("E29").Value = '  +2.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0725'
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.11'
$ws.Range("E32").Value = '  -3.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.55'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("E36").Value = '  -5.05%  '
$ws.Range("E37").Value = '  -0.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.00'
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.59'
$ws.Range("E39").Value = '  -2.69%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.26'
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.97'
$ws.Range("E41").Value = '  -1.33%  '
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.375'
$ws.Range("E42").Value = '  -3.79%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.61'
$ws.Range("E43").Value = '  -2.15%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '281.03'
$ws.Range("E44").Value = '  -1.03%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.08'
$ws.Range("E45").Value = '  -5.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0931'
$ws.Range("E46").Value = '  -1.36%  '
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.93'
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("E51").Value = '  +5.07%  '
